# Sheets update: refresh cached Universalis market-price + Leve profit figures
# (currentAveragePrice / *NQ / *HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ,
#  columns H-N). This workbook stores these as plain numbers (no formulas), so the
# scheduled runner simply overwrites each affected cell with the newly fetched value.

$wb = $excel.ActiveWorkbook

# ============================ ALC ============================
$ws = $wb.Worksheets.Item("ALC")

# Row 15: Morning Glass of Ether (Ether)
$ws.Range("H15").Value = 1555.16
$ws.Range("I15").Value = 1555.16
$ws.Range("K15").Value = 4665.48
$ws.Range("M15").Value = -4496.48

# Row 34: Sophomore Slump (Goatskin Grimoire)
$ws.Range("H34").Value = 2848
$ws.Range("I34").Value = 2848
$ws.Range("K34").Value = 2848
$ws.Range("M34").Value = -2645

# Row 36: You Put Your Left Hand In (Engraved Goatskin Grimoire)
$ws.Range("H36").Value = 2848
$ws.Range("I36").Value = 2848
$ws.Range("K36").Value = 2848
$ws.Range("M36").Value = -2133

# Row 62: The Mustache Suits Him (Enchanted Mythrite Ink)
$ws.Range("H62").Value = 4126.4287
$ws.Range("I62").Value = 2157.8572
$ws.Range("J62").Value = 4618.5713
$ws.Range("K62").Value = 2157.8572
$ws.Range("L62").Value = 4618.5713
$ws.Range("M62").Value = -1533.8572
$ws.Range("N62").Value = -5866.5713

# Row 65: Forgery of Convenience (L) (Enchanted Mythrite Ink)
$ws.Range("H65").Value = 4126.4287
$ws.Range("I65").Value = 2157.8572
$ws.Range("J65").Value = 4618.5713
$ws.Range("K65").Value = 10789.286
$ws.Range("L65").Value = 23092.8565
$ws.Range("M65").Value = -7669.286
$ws.Range("N65").Value = -29332.8565

# Row 100: Asking for a Friend (Beetle Glue)
$ws.Range("H100").Value = 1540.2632
$ws.Range("I100").Value = 1376.0714
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1376.0714
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -835.0714
$ws.Range("N100").Value = -3082

# ============================ ARM ============================
$ws = $wb.Worksheets.Item("ARM")

# Row 32: Ingot We Trust (Steel Ingot)
$ws.Range("H32").Value = 5348.2354
$ws.Range("I32").Value = 4400.0127
$ws.Range("K32").Value = 4400.0127
$ws.Range("M32").Value = -4113.0127

# Row 131: Additions to the Armoire (Chondrite Top of Maiming)
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# ============================ CRP ============================
$ws = $wb.Worksheets.Item("CRP")

# Row 7: Gridania's Got Talent (Maple Lumber)
$ws.Range("H7").Value = 18.272728

# Row 58: You Do the Heavy Lifting (Mahogany Lumber)
$ws.Range("H58").Value = 3647.1777
$ws.Range("I58").Value = 1234.9375
$ws.Range("J58").Value = 4978.069
$ws.Range("K58").Value = 1234.9375
$ws.Range("L58").Value = 4978.069
$ws.Range("M58").Value = -1031.9375
$ws.Range("N58").Value = -5384.069

# Row 136: Turali Quality (Dark Mahogany Lumber)
$ws.Range("H136").Value = 3647.1777
$ws.Range("I136").Value = 1234.9375
$ws.Range("J136").Value = 4978.069
$ws.Range("K136").Value = 3704.8125
$ws.Range("L136").Value = 14934.207
$ws.Range("M136").Value = -1154.8125
$ws.Range("N136").Value = -20034.207

# ============================ CUL ============================
$ws = $wb.Worksheets.Item("CUL")

# Row 104: Fits to a Tea (Doman Tea)
$ws.Range("H104").Value = 26
$ws.Range("I104").Value = 26
$ws.Range("K104").Value = 78
$ws.Range("M104").Value = 2543

# Row 109: Cure for What Ails (Purple Carrot Juice)
$ws.Range("H109").Value = 2380
$ws.Range("I109").Value = 1256
$ws.Range("J109").Value = 3182.8572
$ws.Range("K109").Value = 3768
$ws.Range("L109").Value = 9548.571599999999
$ws.Range("M109").Value = -2728
$ws.Range("N109").Value = -11628.5716

# Row 110: His Dark Utensils (Spaghetti al Nero)
$ws.Range("H110").Value = 2170.6
$ws.Range("I110").Value = 1738.25
$ws.Range("J110").Value = 3900
$ws.Range("K110").Value = 5214.75
$ws.Range("L110").Value = 11700
$ws.Range("M110").Value = -1124.75
$ws.Range("N110").Value = -19880

# Row 113: Can't Eat Just One (Night Vinegar)
$ws.Range("H113").Value = 533.7931
$ws.Range("I113").Value = 535.0476
$ws.Range("J113").Value = 530.5
$ws.Range("K113").Value = 1605.1428
$ws.Range("L113").Value = 1591.5
$ws.Range("M113").Value = 564.8571999999999
$ws.Range("N113").Value = -5931.5

# Row 116: On a Full Stomach (Sausage Links)
$ws.Range("H116").Value = 2230
$ws.Range("I116").Value = 100
$ws.Range("J116").Value = 2466.6667
$ws.Range("K116").Value = 300
$ws.Range("L116").Value = 7400.000100000001
$ws.Range("M116").Value = 3142
$ws.Range("N116").Value = -14284.0001

# Row 118: Teetotally (Masala Chai)
$ws.Range("H118").Value = 1542.8334
$ws.Range("I118").Value = 589.25
$ws.Range("K118").Value = 1767.75
$ws.Range("M118").Value = -524.75

# Row 120: A Happy End (Paella)
$ws.Range("H120").Value = 17333.334
$ws.Range("I120").Value = 8000
$ws.Range("K120").Value = 24000
$ws.Range("M120").Value = -19162

# Row 121: A Cookie for Your Troubles (Coffee Biscuit)
$ws.Range("H121").Value = 50000444
$ws.Range("I121").Value = 495.44446
$ws.Range("J121").Value = 500000000
$ws.Range("K121").Value = 1486.33338
$ws.Range("L121").Value = 1500000000
$ws.Range("M121").Value = -176.33338
$ws.Range("N121").Value = -1500002620

# Row 122: Salt of the North (Northern Sea Salt)
$ws.Range("H122").Value = 9804592
$ws.Range("I122").Value = 17544406
$ws.Range("J122").Value = 826.4666999999999
$ws.Range("K122").Value = 157899654
$ws.Range("L122").Value = 7438.2003
$ws.Range("M122").Value = -157897204
$ws.Range("N122").Value = -12338.2003

# Row 123: Topping Up the Pot (Zurek)
$ws.Range("H123").Value = 2655
$ws.Range("I123").Value = 1715
$ws.Range("J123").Value = 3125
$ws.Range("K123").Value = 5145
$ws.Range("L123").Value = 9375
$ws.Range("M123").Value = -2695
$ws.Range("N123").Value = -14275

# Row 125: At Any Temperature (Borscht)
$ws.Range("H125").Value = 2190
$ws.Range("J125").Value = 3000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -18840

# Row 131: The Mountain Steeped (Tsai tou Vounou)
$ws.Range("H131").Value = 2273.013
$ws.Range("J131").Value = 1732.0571
$ws.Range("L131").Value = 5196.1713
$ws.Range("N131").Value = -15276.1713

# ============================ GSM ============================
$ws = $wb.Worksheets.Item("GSM")

# Row 40: A Little Bird Told Me (Malachite Bracelet)
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5302

# ============================ LTW ============================
$ws = $wb.Worksheets.Item("LTW")

# Row 7: Tan Before the Ban (Leather)
$ws.Range("H7").Value = 1131.3793
$ws.Range("I7").Value = 1066.9048
$ws.Range("J7").Value = 1300.625
$ws.Range("K7").Value = 1066.9048
$ws.Range("L7").Value = 1300.625
$ws.Range("M7").Value = -954.9048
$ws.Range("N7").Value = -1524.625

# Row 22: Skin off Their Backs (Aldgoat Leather)
$ws.Range("H22").Value = 446
$ws.Range("J22").Value = 665
$ws.Range("L22").Value = 665
$ws.Range("N22").Value = -1255

# Row 27: Fire and Hide (Aldgoat Leather)
$ws.Range("H27").Value = 446
$ws.Range("J27").Value = 665
$ws.Range("L27").Value = 665
$ws.Range("N27").Value = -879

# Row 126: Battered Books (Saiga Leather)
$ws.Range("H126").Value = 1131.3793
$ws.Range("I126").Value = 1066.9048
$ws.Range("J126").Value = 1300.625
$ws.Range("K126").Value = 3200.7144
$ws.Range("L126").Value = 3901.875
$ws.Range("M126").Value = -730.7143999999998
$ws.Range("N126").Value = -8841.875

# ============================ WVR ============================
$ws = $wb.Worksheets.Item("WVR")

# Row 39: By the Short Hairs (Velveteen Robe)
$ws.Range("H39").Value = 14900
$ws.Range("J39").Value = 14900
$ws.Range("L39").Value = 14900
$ws.Range("N39").Value = -15726

# Row 42: Put on Your Party Pants (Velveteen Gaskins)
$ws.Range("H42").Value = 13166
$ws.Range("J42").Value = 13166
$ws.Range("L42").Value = 13166
$ws.Range("N42").Value = -13922

# Row 43: Walk Softly and Carry a Big Halberd (Velveteen Dress Shoes)
$ws.Range("H43").Value = 11999
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 13998
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 13998
$ws.Range("M43").Value = -9851
$ws.Range("N43").Value = -14296

# Row 96: Skills on Display (Ruby Cotton Cloth)
$ws.Range("H96").Value = 18250
$ws.Range("I96").Value = 1875
$ws.Range("J96").Value = 51000
$ws.Range("K96").Value = 1875
$ws.Range("L96").Value = 51000
$ws.Range("M96").Value = -502
$ws.Range("N96").Value = -53746
